# "Ajuste de Linea pie de pagina"
#
# - Split the single header/footer into even / default(primary) / first-page
#   variants (Word creates header1-3.xml / footer1-3.xml and wires up the
#   headerReference / footerReference entries in sectPr automatically).
# - The pre-existing header/footer content (the logo + "Instrumento de Uso
#   Oficial..." line) becomes the "default" header/footer; the new even- and
#   first-page header/footer start out blank, using the same paragraph
#   styles as the original ones.
# - Move the stray "_GoBack" bookmark from the body into the (now relocated)
#   default footer paragraph.
# - Nudge the footer distance from the page bottom from 0 to 113 twips
#   (5.65 pt).

$d = $word.ActiveDocument
$s = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2, wdHeaderFooterEvenPages = 3

# --- Even-page header/footer: new, blank, same styles as the originals ---
$evenHeader = $s.Headers.Item(3)
$evenHeader.Range.Text = ""
$evenHeader.Range.Style = "Encabezado"

$evenFooter = $s.Footers.Item(3)
$evenFooter.Range.Text = ""
$evenFooter.Range.Style = "Piedepgina"

# --- First-page header/footer: new, blank, same styles as the originals ---
$firstHeader = $s.Headers.Item(2)
$firstHeader.Range.Text = ""
$firstHeader.Range.Style = "Encabezado"

$firstFooter = $s.Footers.Item(2)
$firstFooter.Range.Text = ""
$firstFooter.Range.Style = "Piedepgina"

# --- Relocate the "_GoBack" bookmark from the body paragraph into the
#     (now split-out) default/primary footer, right before its content. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$primaryFooter = $s.Footers.Item(1)
$footerStart = $primaryFooter.Range
$footerStart.Collapse(1)   # wdCollapseStart
$d.Bookmarks.Add("_GoBack", $footerStart)

# --- Footer margin: 0 -> 113 twips (113/20 pt) ---
$s.PageSetup.FooterDistance = 113 / 20.0
